$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-14 Tuesday" "2025-10-15 Wednesday"

Replace-Text "986÷4=" "323÷4="
Replace-Text "231÷7=" "997÷4="
Replace-Text "644÷7=" "850÷7="
Replace-Text "497÷2=" "833÷8="
Replace-Text "305÷5=" "385÷4="
Replace-Text "830÷6=" "851÷5="
Replace-Text "385÷6=" "395÷3="
Replace-Text "349÷4=" "648÷6="
Replace-Text "871÷3=" "209÷6="
Replace-Text "741÷8=" "468÷9="
Replace-Text "640÷5=" "548÷2="
Replace-Text "606÷3=" "739÷9="
Replace-Text "760÷3=" "149÷5="
Replace-Text "920÷2=" "643÷6="
Replace-Text "849÷3=" "524÷9="
Replace-Text "979÷2=" "349÷6="
Replace-Text "648÷5=" "740÷7="
Replace-Text "881÷9=" "371÷3="
Replace-Text "233÷9=" "551÷5="
Replace-Text "482÷7=" "629÷7="
Replace-Text "306÷2=" "292÷5="
Replace-Text "449÷9=" "424÷8="
Replace-Text "190÷3=" "819÷4="
Replace-Text "569÷6=" "175÷7="
Replace-Text "737÷9=" "565÷9="

Write-Output "Done"
